# HEC2 IN-2 script, only bender - update K4/L4 values, widen column L, move selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the bender-only script values for HEC2 IN-2 (row 4)
$ws.Range("K4").Value = "30.8(22)"
$ws.Range("L4").Value = "2.2(0.18-0.32)"

# Widen column L (12) so the new longer value fits; XML width ends up at 15
$ws.Columns.Item(12).ColumnWidth = 14.166666666666666

# Move the active selection to K5
$ws.Range("K5").Select()
